# Updates the cryptos price/volume table (GitHub Actions refresh).
# Price-looking values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (e.g. "1.001", "5.140") instead of coercing
# them to numbers and dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.574.05'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '1.870.46'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''247.35'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''0.4736'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').Value = '''0.2907'
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').Value = '''0.06468'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '''22.06'
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('D11').Value = '''0.07715'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('D12').Value = '''0.7384'
$ws.Range('E12').Value = '  +4.06%  '
$ws.Range('D13').Value = '''96.46'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('D14').Value = '1.869.54'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '''5.140'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '''272.50'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '30.551.55'
$ws.Range('D18').Value = '''13.30'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '''0.9996'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '''0.000007494'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').Value = '2.118.14'
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = '''5.234'
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('D24').Value = '''6.172'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '''9.206'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').Value = '''163.30'
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('D27').Value = '''18.74'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').Value = '''1.910'
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('D29').Value = '''0.09968'
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').Value = '''1.345'
$ws.Range('E30').Value = '  -2.40%  '
$ws.Range('D31').Value = '''1.507'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').Value = '''4.283'
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').Value = '''4.101'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('D34').Value = '''0.04780'
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('D35').Value = '''1.117'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').Value = '''0.6949'
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').Value = '''1.000'
$ws.Range('D38').Value = '''2.718'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').Value = '''0.01845'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('D41').Value = '''6.180'
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('D42').Value = '''73.15'
$ws.Range('E42').Value = '  +4.19%  '
$ws.Range('D43').Value = '''1.965'
$ws.Range('E43').Value = '  +3.19%  '
$ws.Range('D44').Value = '''0.4169'
$ws.Range('E44').Value = '  +1.76%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '''0.8327'
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').Value = '''102.54'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('D48').Value = '''9.302'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '''6.946'
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '''919.11'
$ws.Range('E51').Value = '  +0.16%  '
